$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New "Arrays" problems: row -> (Sr.No, Title, Url)
# Rows 18 and 19 already existed (only column C was blank); rows 20-29 are
# brand new rows that need Sr.No/Topic, and rows 20-24 also get a linked
# Problem title, while 25-29 stay without one (25-28 empty C cell, 29 no
# C cell at all).
# ---------------------------------------------------------------------------
$rowInfo = @(
  @{ Row = 18; Sr = 17; HasTopic = $false; Title = "Difficulty Rating Order";    Url = "https://www.geeksforgeeks.org/problems/difficulty-rating-order/1" },
  @{ Row = 19; Sr = 18; HasTopic = $false; Title = "Cost of Groceries";          Url = "https://www.geeksforgeeks.org/problems/cost-of-groceries/1" },
  @{ Row = 20; Sr = 19; HasTopic = $true;  Title = "Running Comparison";         Url = "https://www.geeksforgeeks.org/problems/running-comparison/1" },
  @{ Row = 21; Sr = 20; HasTopic = $true;  Title = "Codechef Streak";            Url = "https://www.geeksforgeeks.org/problems/codechef-streak/1" },
  @{ Row = 22; Sr = 21; HasTopic = $true;  Title = "Left Rotate an Array by One";Url = "https://practice.geeksforgeeks.org/problems/left-rotate-an-array-by-one3912/1" },
  @{ Row = 23; Sr = 22; HasTopic = $true;  Title = "Rotate Array by D Steps";    Url = "https://practice.geeksforgeeks.org/problems/rotate-array-by-n-elements-1587115621/1" },
  @{ Row = 24; Sr = 23; HasTopic = $true;  Title = "Leader in an Array";         Url = "https://practice.geeksforgeeks.org/problems/leaders-in-an-array-1587115620/1" },
  @{ Row = 25; Sr = 24; HasTopic = $true;  Title = $null; Url = $null },
  @{ Row = 26; Sr = 25; HasTopic = $true;  Title = $null; Url = $null },
  @{ Row = 27; Sr = 26; HasTopic = $true;  Title = $null; Url = $null },
  @{ Row = 28; Sr = 27; HasTopic = $true;  Title = $null; Url = $null },
  @{ Row = 29; Sr = 28; HasTopic = $true;  Title = $null; Url = $null }
)

foreach ($info in $rowInfo) {
  $row = $info.Row

  # --- Column A: Sr. No. ------------------------------------------------
  if ($info.HasTopic) {
    # Brand-new row: copy formatting from A19 (style s="1") then set value.
    $ws.Range("A19").Copy()
    $ws.Range("A$row").PasteSpecial(-4122) | Out-Null
  }
  $ws.Range("A$row").Value = $info.Sr

  # --- Column B: Topic ("Arrays") ---------------------------------------
  if ($info.HasTopic) {
    $ws.Range("B19").Copy()
    $ws.Range("B$row").PasteSpecial(-4122) | Out-Null
    $ws.Range("B$row").Value = $ws.Range("B19").Value()
  }
}

# --- Column C: Problem (title + hyperlink) for rows 18-24 -----------------
foreach ($info in $rowInfo) {
  if ($null -ne $info.Title) {
    $row = $info.Row

    # Normalize starting formatting to the plain centered style (same as
    # C18/C19 already had before this edit, i.e. the same style as C4)
    # *before* adding the hyperlink. This keeps every Hyperlinks.Add() call
    # starting from an identical base style, so Excel only ever has to mint
    # a single extra transient style behind the scenes instead of one per
    # distinct starting style.
    $ws.Range("C4").Copy()
    $ws.Range("C$row").PasteSpecial(-4122) | Out-Null
    $ws.Range("C$row").ClearContents()

    $ws.Range("C$row").Value = $info.Title
    $ws.Hyperlinks.Add($ws.Range("C$row"), $info.Url) | Out-Null
    # Re-apply the formatting used by the existing hyperlink cells so the
    # cell keeps the same style index as the rest of column C, rather than
    # whatever brand-new style Hyperlinks.Add() just created.
    $ws.Range("C17").Copy()
    $ws.Range("C$row").PasteSpecial(-4122) | Out-Null
  }
}

# --- Column C: empty but styled cells for rows 25-28 -----------------------
# Use C4's formatting (plain centered style, s="1") as the source - this is
# the same style the never-linked C18/C19 cells had before this edit, i.e.
# the "Problem" cell is blank and not yet styled as a hyperlink.
foreach ($row in 25..28) {
  $ws.Range("C4").Copy()
  $ws.Range("C$row").PasteSpecial(-4122) | Out-Null
  $ws.Range("C$row").ClearContents()
}
# Row 29 gets no Problem cell at all - nothing to do for column C there.

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Sheet view: drop the "topLeftCell=A10" scroll position and move the active
# selection down to the new last data cell, C25.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("C25").Select() | Out-Null
